# Adds three new bulleted list items (reusing the same list / numbering
# as the preceding bullets), followed by an empty paragraph and a closing
# plain-text paragraph, right after the last existing paragraph in the
# document ("InfoVisitas - Asesor: ...").

$d = $word.ActiveDocument

$paragraphs = $d.Paragraphs
$count = $paragraphs.Count
$lastParagraph = $paragraphs.Item($count)
$lastEnd = $lastParagraph.Range.End

$bullet1 = "Compra – Inmueble: Una compra puede tener uno o más inmuebles, y un inmueble puede estar contenido en una compra (1:N)"
$bullet2 = "Alquiler – Inmueble: Un alquiler puede ofrecer uno o más inmuebles, y un inmueble puede ser ofrecido por un alquiler (1:N)"
$bullet3 = "Inmueble – Propietario: Un inmueble posee uno o más propietarios y un propietario posee uno o más inmuebles (N:M)"
$closingText = "Procedemos a crear el archivo de workbench con la transformación del modelo entidad-relación"

$listParaXml = @'
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
  </w:pPr>
  <w:r>
    <w:t>{0}</w:t>
  </w:r>
</w:p>
'@

$bulletsXml = ([string]::Format($listParaXml, $bullet1)) + `
              ([string]::Format($listParaXml, $bullet2)) + `
              ([string]::Format($listParaXml, $bullet3))

$tailXml = "<w:p/><w:p><w:r><w:t>$closingText</w:t></w:r></w:p>"

$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + `
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
$bulletsXml + $tailXml + `
'</w:body></w:document>' + `
'</pkg:xmlData></pkg:part></pkg:package>'

# Insert right before the very last paragraph mark of the document body
# (i.e. immediately after the last existing paragraph's text), so the
# existing paragraph is left completely untouched.
$insertionPoint = $d.Range($lastEnd - 1, $lastEnd - 1)
$null = $insertionPoint.InsertXML($payload)
